$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename (A1:D1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the Spanish connector words ("de", "del", "la", "los", "el", "y") ---
# --- in municipality / state names (not the first word of the cell) ---
$ws.Range("B7").Value = "Pabellón De Arteaga"
$ws.Range("B8").Value = "Rincón De Romos"
$ws.Range("B9").Value = "San Francisco De Los Romo"
$ws.Range("B28").Value = "Bejucal De Ocampo"
$ws.Range("B33").Value = "Comitán De Domínguez"
$ws.Range("B45").Value = "Mazapa De Madero"
$ws.Range("B48").Value = "Ocozocoautla De Espinosa"
$ws.Range("B72").Value = "Coyame Del Sotol"
$ws.Range("B77").Value = "Guadalupe Y Calvo"
$ws.Range("B80").Value = "Hidalgo Del Parral"
$ws.Range("B88").Value = "Valle De Zaragoza"
$ws.Range("B100").Value = "San Juan De Sabinas"
$ws.Range("B110").Value = "Villa De Álvarez"
$ws.Range("A112").Value = "Ciudad De México"
$ws.Range("B116").Value = "Cuajimalpa De Morelos"
$ws.Range("B130").Value = "Coneto De Comonfort"
$ws.Range("B144").Value = "Nombre De Dios"
$ws.Range("B153").Value = "San Juan De Guadalupe"
$ws.Range("B154").Value = "San Juan Del Río"
$ws.Range("B155").Value = "San Luis Del Cordero"
$ws.Range("A164").Value = "Estado De México"
$ws.Range("B164").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B170").Value = "Atizapán De Zaragoza"
$ws.Range("B174").Value = "Coacalco De Berriozábal"
$ws.Range("B178").Value = "Ecatepec De Morelos"
$ws.Range("B181").Value = "Ixtapan De La Sal"
$ws.Range("B188").Value = "Naucalpan De Juárez"
$ws.Range("B194").Value = "San Felipe Del Progreso"
$ws.Range("B201").Value = "Tenango Del Valle"
$ws.Range("B203").Value = "Tlalnepantla De Baz"
$ws.Range("B207").Value = "Valle De Bravo"
$ws.Range("B208").Value = "Valle De Chalco Solidaridad"
$ws.Range("B209").Value = "Villa De Allende"
$ws.Range("B218").Value = "San Miguel De Allende"
$ws.Range("B219").Value = "Apaseo El Alto"
$ws.Range("B220").Value = "Apaseo El Grande"
$ws.Range("B228").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B232").Value = "Jaral Del Progreso"
$ws.Range("B239").Value = "Purísima Del Rincón"
$ws.Range("B243").Value = "San Diego De La Unión"
$ws.Range("B245").Value = "San Francisco Del Rincón"
$ws.Range("B247").Value = "San Luis De La Paz"
$ws.Range("B248").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B249").Value = "Silao De La Victoria"
$ws.Range("B252").Value = "Valle De Santiago"
$ws.Range("B258").Value = "Acapulco De Juárez"
$ws.Range("B260").Value = "Ajuchitlán Del Progreso"
$ws.Range("B261").Value = "Alcozauca De Guerrero"
$ws.Range("B264").Value = "Atoyac De Álvarez"
$ws.Range("B265").Value = "Ayutla De Los Libres"
$ws.Range("B268").Value = "Chilapa De Álvarez"
$ws.Range("B269").Value = "Chilpancingo De Los Bravo"
$ws.Range("B273").Value = "Coyuca De Benítez"
$ws.Range("B274").Value = "Coyuca De Catalán"
$ws.Range("B277").Value = "Cuetzala Del Progreso"
$ws.Range("B278").Value = "Cutzamala De Pinzón"
$ws.Range("B284").Value = "Huitzuco De Los Figueroa"
$ws.Range("B285").Value = "Iguala De La Independencia"
$ws.Range("B287").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B288").Value = "Zihuatanejo De Azueta"
$ws.Range("B302").Value = "Taxco De Alarcón"
$ws.Range("B304").Value = "Técpan De Galeana"
$ws.Range("B306").Value = "Tepecoacuilco De Trujano"
$ws.Range("B310").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B322").Value = "Cuautepec De Hinojosa"
$ws.Range("B325").Value = "Huejutla De Reyes"
$ws.Range("B328").Value = "Jacala De Ledezma"
$ws.Range("B332").Value = "Mineral De La Reforma"
$ws.Range("B333").Value = "Mixquiahuala De Juárez"
$ws.Range("B335").Value = "Pachuca De Soto"
$ws.Range("B337").Value = "Progreso De Obregón"
$ws.Range("B338").Value = "Santiago De Anaya"
$ws.Range("B343").Value = "Tepehuacán De Guerrero"
$ws.Range("B344").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B348").Value = "Tula De Allende"
$ws.Range("B349").Value = "Tulancingo De Bravo"
$ws.Range("B350").Value = "Zacualtipán De Ángeles"
$ws.Range("B357").Value = "Atotonilco El Alto"
$ws.Range("B366").Value = "Encarnación De Díaz"
$ws.Range("B370").Value = "Huejuquilla El Alto"
$ws.Range("B371").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B374").Value = "Jilotlán De Los Dolores"
$ws.Range("B377").Value = "Lagos De Moreno"
$ws.Range("B381").Value = "Ojuelos De Jalisco"
$ws.Range("B384").Value = "San Diego De Alejandría"
$ws.Range("B386").Value = "San Martín De Bolaños"
$ws.Range("B388").Value = "Santa María De Los Ángeles"
$ws.Range("B389").Value = "Santa María Del Oro"
$ws.Range("B392").Value = "Talpa De Allende"
$ws.Range("B393").Value = "Tamazula De Gordiano"
$ws.Range("B394").Value = "Techaluta De Montenegro"
$ws.Range("B396").Value = "Teocuitatlán De Corona"
$ws.Range("B397").Value = "Tepatitlán De Morelos"
$ws.Range("B399").Value = "Tizapán El Alto"
$ws.Range("B400").Value = "Tlajomulco De Zúñiga"
$ws.Range("B405").Value = "Unión De San Antonio"
$ws.Range("B406").Value = "Unión De Tula"
$ws.Range("B409").Value = "Zapotlán El Grande"
$ws.Range("B468").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B483").Value = "Coatlán Del Río"
$ws.Range("B490").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B493").Value = "Puente De Ixtla"
$ws.Range("B504").Value = "Amatlán De Cañas"
$ws.Range("B506").Value = "Santa María Del Oro"
$ws.Range("B520").Value = "Mier Y Noriega"
$ws.Range("B526").Value = "San Nicolás De Los Garza"
$ws.Range("B529").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B532").Value = "Constancia Del Rosario"
$ws.Range("B534").Value = "Cuilápam De Guerrero"
$ws.Range("B535").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B536").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B537").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B539").Value = "Ixtlán De Juárez"
$ws.Range("B540").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B542").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B544").Value = "Oaxaca De Juárez"
$ws.Range("B545").Value = "Ocotlán De Morelos"
$ws.Range("B547").Value = "Putla Villa De Guerrero"
$ws.Range("B551").Value = "San Antonio De La Cal"
$ws.Range("B605").Value = "Tataltepec De Valdés"
$ws.Range("B606").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B607").Value = "Tlacolula De Matamoros"
$ws.Range("B608").Value = "Villa De Etla"
$ws.Range("B609").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B610").Value = "Villa De Tututepec"
$ws.Range("B611").Value = "Villa Sola De Vega"
$ws.Range("B621").Value = "Ayotoxco De Guerrero"
$ws.Range("B638").Value = "Ixcamilpa De Guerrero"
$ws.Range("B641").Value = "Izúcar De Matamoros"
$ws.Range("B647").Value = "Palmar De Bravo"
$ws.Range("B655").Value = "San Nicolás De Los Ranchos"
$ws.Range("B664").Value = "Tetela De Ocampo"
$ws.Range("B680").Value = "Amealco De Bonfil"
$ws.Range("B685").Value = "Jalpan De Serra"
$ws.Range("B686").Value = "Landa De Matamoros"
$ws.Range("B689").Value = "Pinal De Amoles"
$ws.Range("B691").Value = "San Juan Del Río"
$ws.Range("B699").Value = "Axtla De Terrazas"
$ws.Range("B704").Value = "Ciudad Del Maíz"
$ws.Range("B711").Value = "Mexquitic De Carmona"
$ws.Range("B716").Value = "San Ciro De Acosta"
$ws.Range("B721").Value = "Santa María Del Río"
$ws.Range("B728").Value = "Tanquián De Escobedo"
$ws.Range("B731").Value = "Villa De Arista"
$ws.Range("B732").Value = "Villa De Arriaga"
$ws.Range("B733").Value = "Villa De Guadalupe"
$ws.Range("B734").Value = "Villa De Ramos"
$ws.Range("B735").Value = "Villa De Reyes"
$ws.Range("B797").Value = "Soto La Marina"
$ws.Range("B804").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B811").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B833").Value = "Castillo De Teayo"
$ws.Range("B834").Value = "Cazones De Herrera"
$ws.Range("B838").Value = "Chinampa De Gorostiza"
$ws.Range("B844").Value = "Cosamaloapan De Carpio"
$ws.Range("B855").Value = "Hueyapan De Ocampo"
$ws.Range("B856").Value = "Ignacio De La Llave"
$ws.Range("B858").Value = "Ixhuatlán Del Café"
$ws.Range("B859").Value = "Ixhuatlán Del Sureste"
$ws.Range("B864").Value = "Juchique De Ferrer"
$ws.Range("B866").Value = "Lerdo De Tejada"
$ws.Range("B867").Value = "Martínez De La Torre"
$ws.Range("B879").Value = "Poza Rica De Hidalgo"
$ws.Range("B885").Value = "Sayula De Alemán"
$ws.Range("B901").Value = "Vega De Alatorre"
$ws.Range("B915").Value = "Cañitas De Felipe Pescador"
$ws.Range("B916").Value = "Concepción Del Oro"
$ws.Range("B923").Value = "Jiménez Del Teul"
$ws.Range("B929").Value = "Noria De Ángeles"
$ws.Range("B939").Value = "Teúl De González Ortega"
$ws.Range("B940").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B943").Value = "Villa De Cos"

# --- Remove trailing footer/metadata rows (951:955) and shrink used range ---
$ws.Range("A951:D955").EntireRow.Delete()

